# Auto-generated files on 2025-10-27
# Updates the stock-name table (A2:C21) on the active worksheet to reflect
# the latest "Hot Stock Top 20" rankings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("江波龙", "闻泰科技", "神开股份"),
    @("平潭发展", "工业富联", "华建集团"),
    @("安泰科技", "江波龙", "江波龙"),
    @("大为股份", "神开股份", "胜宏科技"),
    @("神开股份", "晶瑞电材", "工业富联"),
    @("工业富联", "平潭发展", "平潭发展"),
    @("闻泰科技", "安泰科技", "闻泰科技"),
    @("晶瑞电材", "大为股份", "安泰科技"),
    @("达华智能", "方正科技", "兆易创新"),
    @("胜宏科技", "兆易创新", "达华智能"),
    @("兆易创新", "达华智能", "向日葵"),
    @("盈新发展", "大有能源", "盈新发展"),
    @("中国核建", "东方财富", "恒宝股份"),
    @("大有能源", "胜宏科技", "中信重工"),
    @("恒宝股份", "三花智控", "步步高"),
    @("三花智控", "恒宝股份", "德明利"),
    @("德明利", "新易盛", "三花智控"),
    @("方正科技", "中国核建", "世龙实业"),
    @("华建集团", "郑州煤电", "大为股份"),
    @("中际旭创", "中国卫星", "首开股份")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}
